$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "Updating capital structure database..."
$ws.Rows.Item(4).Delete()

# Row 2
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1"
$ws.Range("B2").Style = "Normal"
$ws.Range("D2").Value = 0.112
$ws.Range("E2").Value = 0.0836
$ws.Range("G2").Value = 0.3314814814814815
$ws.Range("H2").Value = 0.3314814814814815
$ws.Range("I2").Value = 0.3393518518518518
$ws.Range("J2").Value = 0.2428425390708628
$ws.Range("K2").Value = 6.14
$ws.Range("L2").Value = 0.2842592592592592
$ws.Range("M2").Value = 5.87
$ws.Range("N2").Value = 0.04811475409836066
$ws.Range("O2").Value = 0.9560260586319219
$ws.Range("P2").Value = 5.87
$ws.Range("Q2").Value = 0.04811475409836066
$ws.Range("R2").Value = 0.9560260586319219
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 2.77
$ws.Range("V2").Value = 0.02270491803278689
$ws.Range("W2").Value = 0.1035413153456998
$ws.Range("X2").Value = 0.0193322958805533
$ws.Range("Y2").Value = 0.08420901946514653
$ws.Range("Z2").Value = 0.6607525237075559
$ws.Range("AA2").Value = 0.1604588205546233
$ws.Range("AB2").Value = 0.01941958683860227
$ws.Range("AC2").Value = 0.1410392337160211
$ws.Range("AD2").Value = 6.69
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 6.69
$ws.Range("AG2").Value = 3.92
$ws.Range("AH2").Value = 0.0519853912502914
$ws.Range("AI2").Value = 0.1055371509701846
$ws.Range("AJ2").Value = 0.03113087674714105
$ws.Range("AK2").Value = 0.06466512702078522
$ws.Range("AL2").Value = 0.585
$ws.Range("AM2").Value = 0.585
$ws.Range("AN2").Value = 0.8654592496765847
$ws.Range("AO2").Value = 12.52991452991453
$ws.Range("AP2").Value = 0.5071151358344114
$ws.Range("AQ2").Value = 12.52991452991453

# Row 3
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "Bolsa de Valores de Lima S.A.A. (BVL:BVLAC1)"
$ws.Range("B3").Style = "Normal"
$ws.Range("D3").Value = 0.112
$ws.Range("E3").Value = 0.0836
$ws.Range("G3").Value = 0.3314814814814815
$ws.Range("H3").Value = 0.3314814814814815
$ws.Range("I3").Value = 0.3393518518518518
$ws.Range("J3").Value = 0.2428425390708628
$ws.Range("K3").Value = 6.14
$ws.Range("L3").Value = 0.2842592592592592
$ws.Range("M3").Value = 5.87
$ws.Range("N3").Value = 0.04811475409836066
$ws.Range("O3").Value = 0.9560260586319219
$ws.Range("P3").Value = 5.87
$ws.Range("Q3").Value = 0.04811475409836066
$ws.Range("R3").Value = 0.9560260586319219
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 2.77
$ws.Range("V3").Value = 0.02270491803278689
$ws.Range("W3").Value = 0.1035413153456998
$ws.Range("X3").Value = 0.0193322958805533
$ws.Range("Y3").Value = 0.08420901946514653
$ws.Range("Z3").Value = 0.6607525237075559
$ws.Range("AA3").Value = 0.1604588205546233
$ws.Range("AB3").Value = 0.01941958683860227
$ws.Range("AC3").Value = 0.1410392337160211
$ws.Range("AD3").Value = 6.69
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 6.69
$ws.Range("AG3").Value = 3.92
$ws.Range("AH3").Value = 0.0519853912502914
$ws.Range("AI3").Value = 0.1055371509701846
$ws.Range("AJ3").Value = 0.03113087674714105
$ws.Range("AK3").Value = 0.06466512702078522
$ws.Range("AL3").Value = 0.585
$ws.Range("AM3").Value = 0.585
$ws.Range("AN3").Value = 0.8654592496765847
$ws.Range("AO3").Value = 12.52991452991453
$ws.Range("AP3").Value = 0.5071151358344114
$ws.Range("AQ3").Value = 12.52991452991453
